{"js": "// Removed a reference to the PTM chapter in the tutorials.\n//\n// Two edits against the \"This chapter will guide you through all these\n// steps, ...\" paragraph and the heading list that follows it:\n//   1. \"separated into six sections, and finally provide an introduction\n//       to post-translational modifications (PTM) oriented studies:\"\n//      becomes\n//      \"separated into five sections:\"\n//   2. The \"1.6   PTM Analysis\" heading paragraph is deleted entirely.\n\nconst body = context.document.body;\n\n// --- Edit 1: \"six\" -> \"five\" in the \"separated into six sections\" phrase.\n// Search for the unique surrounding phrase so the other, unrelated \"six\"\n// earlier in the document (`...divided into six steps:`) is untouched.\nconst sixSearch = body.search(\"separated into six sections\", { matchCase: true });\nsixSearch.load(\"items\");\nawait context.sync();\n\nif (sixSearch.items.length > 0) {\n  sixSearch.items[0].insertText(\"separated into five sections\", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 2: drop the trailing \", and finally provide an introduction to\n// post-translational modifications (PTM) oriented studies\" clause, leaving\n// just the closing colon.\nconst tailSearch = body.search(\n  \", and finally provide an introduction to post-translational modifications (PTM) oriented studies:\",\n  { matchCase: true }\n);\ntailSearch.load(\"items\");\nawait context.sync();\n\nif (tailSearch.items.length > 0) {\n  tailSearch.items[0].insertText(\":\", \"Replace\");\n  await context.sync();\n}\n\n// --- Edit 3: remove the \"1.6   PTM Analysis\" heading paragraph entirely.\nconst headingSearch = body.search(\"PTM Analysis\", { matchCase: true });\nheadingSearch.load(\"items\");\nawait context.sync();\n\nfor (const hit of headingSearch.items) {\n  const para = hit.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n  if (/^\\s*1\\.6\\s*PTM\\s*Analysis\\s*$/.test(para.text)) {\n    para.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Removed a reference to the PTM chapter in the tutorials.\n#\n# 1. \"...separated into six sections, and finally provide an introduction\n#     to post-translational modifications (PTM) oriented studies:\"\n#    becomes\n#    \"...separated into five sections:\"\n# 2. The \"1.6    PTM Analysis\" heading paragraph is deleted entirely.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: \"six\" -> \"five\" in \"separated into six sections\".\n# Scope the search text to the specific phrase so the earlier, unrelated\n# \"...divided into six steps:\" sentence is left untouched.\n$rng1 = $d.Content\n$rng1.Find.Execute(\"separated into six sections\", $false, $false, $false, $false, $false, $true, 1, $false, \"separated into five sections\", 2)\n\n# --- Edit 2: drop the \", and finally provide an introduction to\n# post-translational modifications (PTM) oriented studies\" clause, leaving\n# just the trailing colon.\n$rng2 = $d.Content\n$rng2.Find.Execute(\", and finally provide an introduction to post-translational modifications (PTM) oriented studies:\", $false, $false, $false, $false, $false, $true, 1, $false, \":\", 2)\n\n# --- Edit 3: remove the \"1.6   PTM Analysis\" heading paragraph entirely,\n# including its paragraph mark, leaving the rest of the heading list intact.\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -match \"PTM Analysis\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
